$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last (18th) data row entirely so the table shrinks from 18 to 17 players
$ws.Range("A19:C19").ClearContents()

# Rewrite the roster table (players reshuffled, "Josh Okogie" and "Bilal Coulibaly" dropped,
# "Andrew Nembhard" added) in the new order shown in the updated file.
$ws.Range("A2").Value = 'D''Angelo Russell'
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Brooklyn Nets"

$ws.Range("A3").Value = "Devin Booker"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Phoenix Suns"

$ws.Range("A4").Value = "Trae Young"
$ws.Range("B4").Value = "PG"
$ws.Range("C4").Value = "Atlanta Hawks"

$ws.Range("A5").Value = "Jalen Brunson"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "New York Knicks"

$ws.Range("A6").Value = "LeBron James"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Los Angeles Lakers"

$ws.Range("A7").Value = "Christian Braun"
$ws.Range("B7").Value = "SG,SF"
$ws.Range("C7").Value = "Denver Nuggets"

$ws.Range("A8").Value = "Myles Turner"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "Indiana Pacers"

$ws.Range("A9").Value = "Walker Kessler"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Utah Jazz"

$ws.Range("A10").Value = "Scoot Henderson"
$ws.Range("B10").Value = "PG"
$ws.Range("C10").Value = "Portland Trail Blazers"

$ws.Range("A11").Value = "Andrew Nembhard"
$ws.Range("B11").Value = "PG,SG"
$ws.Range("C11").Value = "Indiana Pacers"

$ws.Range("A12").Value = "Jalen Williams"
$ws.Range("B12").Value = "SG,SF,PF,C"
$ws.Range("C12").Value = "Oklahoma City Thunder"

$ws.Range("A13").Value = "Norman Powell"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "LA Clippers"

$ws.Range("A14").Value = "Immanuel Quickley"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Toronto Raptors"

$ws.Range("A15").Value = "Kawhi Leonard"
$ws.Range("B15").Value = "SG,SF,PF"
$ws.Range("C15").Value = "LA Clippers"

$ws.Range("A16").Value = "Desmond Bane"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Memphis Grizzlies"

$ws.Range("A17").Value = "Brandon Ingram"
$ws.Range("B17").Value = "SG,SF,PF"
$ws.Range("C17").Value = "New Orleans Pelicans"

$ws.Range("A18").Value = "Jimmy Butler"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Miami Heat"
